# Update cryptocurrency price/volume data (and shifted coin rows 9-51) per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '29.022.18'
$ws.Range("E2").Value = '  +0.02%  '

# Row 3: Ethereum
$ws.Range("D3").Value = '1.833.56'
$ws.Range("E3").Value = '  +0.25%  '

# Row 4: TetherUSD
$ws.Range("D4").Value = '''0.9975'
$ws.Range("E4").Value = '  -0.18%  '

# Row 5: BNB
$ws.Range("D5").Value = '''242.53'
$ws.Range("E5").Value = '  +0.54%  '

# Row 6: XRP
$ws.Range("D6").Value = '''0.6277'
$ws.Range("E6").Value = '  -3.92%  '

# Row 7: USDC
$ws.Range("D7").Value = '''1.0000'
$ws.Range("E7").Value = '  -0.05%  '

# Row 8: Dogecoin
$ws.Range("D8").Value = '''0.07611'
$ws.Range("E8").Value = '  +3.75%  '

# Row 9: Cardano
$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D9").Value = '''0.2924'
$ws.Range("E9").Value = '  -0.30%  '

# Row 10: Solana
$ws.Range("B10").Value = 'Solana'
$ws.Range("C10").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D10").Value = '''22.57'
$ws.Range("E10").Value = '  -1.50%  '

# Row 11: TRON
$ws.Range("B11").Value = 'TRON'
$ws.Range("C11").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D11").Value = '''0.07715'
$ws.Range("E11").Value = '  +0.61%  '

# Row 12: WrappedEther
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.833.67'
$ws.Range("E12").Value = '  +0.17%  '

# Row 13: Polkadot
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = '''4.949'
$ws.Range("E13").Value = '  -0.56%  '

# Row 14: Polygon
$ws.Range("B14").Value = 'Polygon'
$ws.Range("C14").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D14").Value = '''0.6644'
$ws.Range("E14").Value = '  -0.19%  '

# Row 15: ShibaInu
$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D15").Value = '''0.00001019'
$ws.Range("E15").Value = '  +17.79%  '

# Row 16: Litecoin
$ws.Range("B16").Value = 'Litecoin'
$ws.Range("C16").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D16").Value = '''82.84'
$ws.Range("E16").Value = '  +1.16%  '

# Row 17: Uniswap
$ws.Range("B17").Value = 'Uniswap'
$ws.Range("C17").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D17").Value = '''6.053'
$ws.Range("E17").Value = '  -0.33%  '

# Row 18: WrappedBTC
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '29.049.98'
$ws.Range("E18").Value = '  +0.39%  '

# Row 19: BitcoinCash
$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D19").Value = '''226.59'
$ws.Range("E19").Value = '  +1.25%  '

# Row 20: Avalanche
$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D20").Value = '''12.34'
$ws.Range("E20").Value = '  -0.61%  '

# Row 21: Dai
$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").Value = '''0.9986'
$ws.Range("E21").Value = '  -0.16%  '

# Row 22: Chainlink
$ws.Range("B22").Value = 'Chainlink'
$ws.Range("C22").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D22").Value = '''7.196'
$ws.Range("E22").Value = '  +1.38%  '

# Row 23: BinanceUSD
$ws.Range("B23").Value = 'BinanceUSD'
$ws.Range("C23").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D23").Value = '''0.9989'
$ws.Range("E23").Value = '  -0.16%  '

# Row 24: Monero
$ws.Range("B24").Value = 'Monero'
$ws.Range("C24").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D24").Value = '''158.46'
$ws.Range("E24").Value = '  +0.26%  '

# Row 25: Cosmos
$ws.Range("B25").Value = 'Cosmos'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D25").Value = '''8.491'
$ws.Range("E25").Value = '  -0.09%  '

# Row 26: Stellar
$ws.Range("B26").Value = 'Stellar'
$ws.Range("C26").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D26").Value = '''0.1372'
$ws.Range("E26").Value = '  -0.49%  '

# Row 27: EthereumClassic
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").Value = '''17.88'
$ws.Range("E27").Value = '  +0.03%  '

# Row 28: PancakeSwap
$ws.Range("B28").Value = 'PancakeSwap'
$ws.Range("C28").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D28").Value = '''1.490'
$ws.Range("E28").Value = '  -0.85%  '

# Row 29: Filecoin
$ws.Range("B29").Value = 'Filecoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D29").Value = '''4.094'
$ws.Range("E29").Value = '  -0.21%  '

# Row 30: InternetComputer(DFINITY)
$ws.Range("B30").Value = 'InternetComputer(DFINITY)'
$ws.Range("C30").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D30").Value = '''4.015'
$ws.Range("E30").Value = '  +0.10%  '

# Row 31: Toncoin
$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D31").Value = '''1.188'
$ws.Range("E31").Value = '  -1.09%  '

# Row 32: Hedera
$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").Value = '''0.05236'
$ws.Range("E32").Value = '  -2.11%  '

# Row 33: LidoDAOToken
$ws.Range("B33").Value = 'LidoDAOToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D33").Value = '''1.840'
$ws.Range("E33").Value = '  +0.31%  '

# Row 34: ImmutableX
$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D34").Value = '''0.7353'
$ws.Range("E34").Value = '  -1.05%  '

# Row 35: ARBITRUM
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").Value = '''1.138'
$ws.Range("E35").Value = '  -1.21%  '

# Row 36: HuobiToken
$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D36").Value = '''2.700'
$ws.Range("E36").Value = '  +2.20%  '

# Row 37: Maker
$ws.Range("B37").Value = 'Maker'
$ws.Range("C37").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D37").Value = '1.238.52'
$ws.Range("E37").Value = '  -4.78%  '

# Row 38: MXToken
$ws.Range("B38").Value = 'MXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D38").Value = '''2.757'
$ws.Range("E38").Value = '  +0.53%  '

# Row 39: VeChain
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = '''0.01785'
$ws.Range("E39").Value = '  -0.01%  '

# Row 40: FraxShare
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").Value = '''6.370'
$ws.Range("E40").Value = '  -0.03%  '

# Row 41: TrustWalletToken
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").Value = '''0.8962'
$ws.Range("E41").Value = '  +0.10%  '

# Row 42: PaxDollar
$ws.Range("B42").Value = 'PaxDollar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D42").Value = '''1.000'
$ws.Range("E42").Value = '  +0.06%  '

# Row 43: Quant
$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D43").Value = '''102.32'
$ws.Range("E43").Value = '  -0.72%  '

# Row 44: BabyDogeCoin
$ws.Range("B44").Value = 'BabyDogeCoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D44").Value = '''0.00000000125'
$ws.Range("E44").Value = '  +5.40%  '

# Row 45: RocketPoolETH
$ws.Range("B45").Value = 'RocketPoolETH'
$ws.Range("C45").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D45").Value = '1.981.85'
$ws.Range("E45").Value = '  -0.07%  '

# Row 46: Aave
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").Value = '''64.04'
$ws.Range("E46").Value = '  +0.01%  '

# Row 47: Mantle
$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D47").Value = '''0.5100'
$ws.Range("E47").Value = '  -0.77%  '

# Row 48: TheSandbox
$ws.Range("B48").Value = 'TheSandbox'
$ws.Range("C48").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D48").Value = '''0.4038'
$ws.Range("E48").Value = '  +1.41%  '

# Row 49: EnergySwap
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '''8.879'
$ws.Range("E49").Value = '  +1.70%  '

# Row 50: Cronos
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = '''0.05742'
$ws.Range("E50").Value = '  -1.46%  '

# Row 51: Aptos
$ws.Range("B51").Value = 'Aptos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D51").Value = '''6.690'
$ws.Range("E51").Value = '  +0.02%  '

